# Replace the mixed text/number job-shop data in row 1 with plain numeric
# values (matching the format already used by every other row), then move
# the active selection to I26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2,24,3,12,9,17,4,27,0,21,6,25,8,27,7,26,1,30,5,31,11,18,14,16,13,39,10,19,12,26)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}

$ws.Range("I26").Select()
